$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) for column A from row 16 into the new rows 17-19
$ws.Range("A16").Copy($ws.Range("A17"))
$ws.Range("A16").Copy($ws.Range("A18"))
$ws.Range("A16").Copy($ws.Range("A19"))

# --- Column B (Question) first, in row order, so shared-string indices line up ---
$q1 = @'
les seuils ph?
'@
$ws.Range("B17").Value = $q1
$q2 = @'
donne moi la procedure complete en details su le reglage du ph
'@
$ws.Range("B18").Value = $q2
$q3 = @'
donne moi les login et mot de passe pour chaque machiine
'@
$ws.Range("B19").Value = $q3

# --- Column C (Réponse) next ---
$a1 = @'
 
Il n'y a pas de seuils PH spécifiés dans ces données. Les données indiquent que le PH devrait être à 9 et que l'electrode doit être changée pour voir s'il y a du mieux. L'étalonnage doit également être effectué avec SUEZ, et des bidons PH4, PH7 et PH9 doivent être bas pour cela.
'@
$ws.Range("C17").Value = $a1
$a2 = @'
 de la cuve
Procédure pour le réglage du PH de la cuve : 
1. Préparation pour l'étalonnage - Remplir le bidon avec du PH 4, du PH 7 et du PH 9 et placer les fioles pour l'étalonnage en bas. Disposer la pipette demi-eau en bas également.
2. Monter la sonde et procéder à l'étalonnage - Monter la sonde et procéder à l'étalonnage selon les instructions de Suez.
3. Connexion et remise en place du tout - Effectuer la connexion et assurer la remise en place du tout.
4. Attente de la remise en température des cuves à 52°C - Attendre que la température des cuves atteigne les 52°C pour procéder à l'essai avec l'opérateur.
5. Essais avec l'opérateur - Procéder aux essais avec l'opérateur selon les instructions.
'@
$ws.Range("C18").Value = $a2
$a3 = @'

Pour la machine 5295 : Utilisateur: 123  Mot de passe: 321
Pour la machine 51783 : Utilisateur: 752  Mot de passe: 752
Pour la machine 5232 : Aucun login et mot de passe spécifiés
Pour la machine 5252 : Aucun login et mot de passe spécifiés
Pour la machine 5304 : Utilisateur: 752  Mot de passe: 752
'@
$ws.Range("C19").Value = $a3

# --- Column D (Commentaire) last ---
$ws.Range("D17").Value = "bonne reponse"
$c2 = @'
yess
'@
$ws.Range("D18").Value = $c2
$c3 = @'
tres bonne reponse merci
'@
$ws.Range("D19").Value = $c3

# --- Column A (row numbers) and E (Avis) ---
$ws.Range("A17").Value = 15
$ws.Range("A18").Value = 16
$ws.Range("A19").Value = 17

$ws.Range("E17").Value = 5
$ws.Range("E18").Value = 5
$ws.Range("E19").Value = 5
